$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.92" are not
# auto-converted into floating point numbers, matching the original inlineStr text cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "96.470.54"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "3.716.76"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "238.38"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").Value = "1.92"
$ws.Range("E6").Value = "  +6.85%  "
$ws.Range("D7").Value = "655.89"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "0.423"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "3.712.59"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("D12").Value = "45.09"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  +5.89%  "
$ws.Range("D15").Value = "4.409.41"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("D16").Value = "0.0000268"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").Value = "96.310.50"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "9.00"
$ws.Range("E18").Value = "  +16.08%  "
$ws.Range("D19").Value = "3.686.66"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").Value = "19.13"
$ws.Range("E20").Value = "  +4.68%  "
$ws.Range("D21").Value = "12.79"
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("D22").Value = "0.527"
$ws.Range("E22").Value = "  -2.47%  "
$ws.Range("D23").Value = "525.09"
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("D24").Value = "3.50"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "0.0000203"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").Value = "102.64"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").Value = "13.46"
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("D29").Value = "0.169"
$ws.Range("E29").Value = "  -7.30%  "
$ws.Range("D30").Value = "12.49"
$ws.Range("E30").Value = "  +3.65%  "
$ws.Range("D31").Value = "3.07"
$ws.Range("E31").Value = "  +2.47%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "1.91"
$ws.Range("E33").Value = "  +10.79%  "
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("D35").Value = "671.50"
$ws.Range("E35").Value = "  +8.51%  "
$ws.Range("D36").Value = "32.83"
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "0.600"
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("D39").Value = "8.89"
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("D40").Value = "7.12"
$ws.Range("E40").Value = "  +16.11%  "
$ws.Range("E41").Value = "  +4.63%  "
$ws.Range("D42").Value = "40.19"
$ws.Range("E42").Value = "  +22.80%  "
$ws.Range("D43").Value = "0.982"
$ws.Range("E43").Value = "  +5.59%  "
$ws.Range("D44").Value = "1.98"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "0.0461"
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "8.62"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("E51").Value = "  +2.42%  "

# Restore default (General) formatting/style so the cells keep looking like the
# original unstyled cells (no explicit style index), only the text content changed.
$ws.Range("D2:D51").ClearFormats()
